# Auto-generated Excel COM-interop script
# Applies scheduled-runner value updates to the Leve profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 708.0769
$ws.Range("J17").Value = 708.0769
$ws.Range("L17").Value = 2124.2307
$ws.Range("N17").Value = -2460.2307
$ws.Range("H28").Value = 1125.2
$ws.Range("I28").Value = 1110.2142
$ws.Range("J28").Value = 1160.1666
$ws.Range("K28").Value = 1110.2142
$ws.Range("L28").Value = 1160.1666
$ws.Range("M28").Value = -625.2141999999999
$ws.Range("N28").Value = -2130.1666
$ws.Range("H111").Value = 3222.2222
$ws.Range("I111").Value = 2600
$ws.Range("J111").Value = 4000
$ws.Range("K111").Value = 7800
$ws.Range("L111").Value = 12000
$ws.Range("M111").Value = -4733
$ws.Range("N111").Value = -18134
$ws.Range("H125").Value = 15000
$ws.Range("I125").Value = 6000
$ws.Range("J125").Value = 19500
$ws.Range("K125").Value = 54000
$ws.Range("L125").Value = 175500
$ws.Range("M125").Value = -51540
$ws.Range("N125").Value = -180420
$ws.Range("H132").Value = 1276402.2
$ws.Range("I132").Value = 1374306.2
$ws.Range("J132").Value = 3650
$ws.Range("K132").Value = 4122918.6
$ws.Range("L132").Value = 10950
$ws.Range("M132").Value = -4120388.6
$ws.Range("N132").Value = -16010
$ws.Range("H135").Value = 2659.5
$ws.Range("I135").Value = 3122.8572
$ws.Range("J135").Value = 1185.1818
$ws.Range("K135").Value = 28105.7148
$ws.Range("L135").Value = 10666.6362
$ws.Range("M135").Value = -25570.7148
$ws.Range("N135").Value = -15736.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1162.8667
$ws.Range("I2").Value = 1026.6364
$ws.Range("J2").Value = 1537.5
$ws.Range("K2").Value = 1026.6364
$ws.Range("L2").Value = 1537.5
$ws.Range("M2").Value = -913.6364000000001
$ws.Range("N2").Value = -1763.5
$ws.Range("H32").Value = 5302.03
$ws.Range("I32").Value = 3177.854
$ws.Range("J32").Value = 22488.545
$ws.Range("K32").Value = 3177.854
$ws.Range("L32").Value = 22488.545
$ws.Range("M32").Value = -2890.854
$ws.Range("N32").Value = -23062.545
$ws.Range("H102").Value = 1900
$ws.Range("I102").Value = 1900
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1900
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -278
$ws.Range("N102").ClearContents()
$ws.Range("H116").Value = 1162.8667
$ws.Range("I116").Value = 1026.6364
$ws.Range("J116").Value = 1537.5
$ws.Range("K116").Value = 1026.6364
$ws.Range("L116").Value = 1537.5
$ws.Range("M116").Value = 1267.3636
$ws.Range("N116").Value = -6125.5
$ws.Range("H122").Value = 1816
$ws.Range("I122").Value = 1305.92
$ws.Range("J122").Value = 2975.2727
$ws.Range("K122").Value = 3917.76
$ws.Range("L122").Value = 8925.8181
$ws.Range("M122").Value = -1467.76
$ws.Range("N122").Value = -13825.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1162.8667
$ws.Range("I3").Value = 1026.6364
$ws.Range("J3").Value = 1537.5
$ws.Range("K3").Value = 1026.6364
$ws.Range("L3").Value = 1537.5
$ws.Range("M3").Value = -912.6364000000001
$ws.Range("N3").Value = -1765.5
$ws.Range("H64").Value = 816.5
$ws.Range("I64").Value = 613
$ws.Range("J64").Value = 1020
$ws.Range("K64").Value = 613
$ws.Range("L64").Value = 1020
$ws.Range("M64").Value = -388
$ws.Range("N64").Value = -1470
$ws.Range("H67").Value = 816.5
$ws.Range("I67").Value = 613
$ws.Range("J67").Value = 1020
$ws.Range("K67").Value = 613
$ws.Range("L67").Value = 1020
$ws.Range("M67").Value = 167
$ws.Range("N67").Value = -2580
$ws.Range("H88").Value = 27000
$ws.Range("I88").Value = 27000
$ws.Range("K88").Value = 27000
$ws.Range("M88").Value = -26594
$ws.Range("H91").Value = 27000
$ws.Range("I91").Value = 27000
$ws.Range("K91").Value = 27000
$ws.Range("M91").Value = -25596
$ws.Range("H107").Value = 1211.1034
$ws.Range("I107").Value = 1291.4546
$ws.Range("J107").Value = 958.5714
$ws.Range("K107").Value = 1291.4546
$ws.Range("L107").Value = 958.5714
$ws.Range("M107").Value = 628.5454
$ws.Range("N107").Value = -4798.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5939.4614
$ws.Range("I16").Value = 4300
$ws.Range("K16").Value = 4300
$ws.Range("M16").Value = -4013
$ws.Range("H31").Value = 2062.5112
$ws.Range("I31").Value = 1651.8975
$ws.Range("J31").Value = 4731.5
$ws.Range("K31").Value = 1651.8975
$ws.Range("L31").Value = 4731.5
$ws.Range("M31").Value = -1356.8975
$ws.Range("N31").Value = -5321.5
$ws.Range("H34").Value = 2062.5112
$ws.Range("I34").Value = 1651.8975
$ws.Range("J34").Value = 4731.5
$ws.Range("K34").Value = 1651.8975
$ws.Range("L34").Value = 4731.5
$ws.Range("M34").Value = -1449.8975
$ws.Range("N34").Value = -5135.5
$ws.Range("H99").Value = 2836.4546
$ws.Range("I99").Value = 1885.8572
$ws.Range("J99").Value = 4500
$ws.Range("K99").Value = 1885.8572
$ws.Range("L99").Value = 4500
$ws.Range("M99").Value = -387.8571999999999
$ws.Range("N99").Value = -7496
$ws.Range("H113").Value = 5939.4614
$ws.Range("I113").Value = 4300
$ws.Range("K113").Value = 4300
$ws.Range("M113").Value = -2130
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H122").Value = 1599.3846
$ws.Range("I122").Value = 1713.1428
$ws.Range("J122").Value = 1466.6666
$ws.Range("K122").Value = 5139.428400000001
$ws.Range("L122").Value = 4399.9998
$ws.Range("M122").Value = -2689.428400000001
$ws.Range("N122").Value = -9299.9998
$ws.Range("H126").Value = 2836.4546
$ws.Range("I126").Value = 1885.8572
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 5657.571599999999
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -3187.571599999999
$ws.Range("N126").Value = -18440

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 50502
$ws.Range("I68").Value = 50502
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 151506
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -150695
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 50502
$ws.Range("I71").Value = 50502
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 454518
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -450462
$ws.Range("N71").ClearContents()
$ws.Range("H80").Value = 2852.0952
$ws.Range("I80").Value = 1451
$ws.Range("J80").Value = 2999.5789
$ws.Range("K80").Value = 4353
$ws.Range("L80").Value = 8998.736699999999
$ws.Range("M80").Value = -3417
$ws.Range("N80").Value = -10870.7367
$ws.Range("H83").Value = 2852.0952
$ws.Range("I83").Value = 1451
$ws.Range("J83").Value = 2999.5789
$ws.Range("K83").Value = 13059
$ws.Range("L83").Value = 26996.2101
$ws.Range("M83").Value = -8379
$ws.Range("N83").Value = -36356.2101
$ws.Range("H86").Value = 500
$ws.Range("J86").Value = 500
$ws.Range("L86").Value = 1500
$ws.Range("N86").Value = -3872
$ws.Range("H89").Value = 500
$ws.Range("J89").Value = 500
$ws.Range("L89").Value = 4500
$ws.Range("N89").Value = -16356
$ws.Range("H92").Value = 1106.25
$ws.Range("I92").Value = 537.5
$ws.Range("K92").Value = 1612.5
$ws.Range("M92").Value = -364.5
$ws.Range("H107").Value = 231.4
$ws.Range("I107").Value = 115.73333
$ws.Range("J107").Value = 404.9
$ws.Range("K107").Value = 347.19999
$ws.Range("L107").Value = 1214.7
$ws.Range("M107").Value = 1572.80001
$ws.Range("N107").Value = -5054.7
$ws.Range("H122").Value = 2381445.8
$ws.Range("I122").Value = 323.2857
$ws.Range("J122").Value = 7143691
$ws.Range("K122").Value = 2909.5713
$ws.Range("L122").Value = 64293219
$ws.Range("M122").Value = -459.5713000000001
$ws.Range("N122").Value = -64298119

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 25689.428
$ws.Range("I26").Value = 20000
$ws.Range("J26").Value = 26637.666
$ws.Range("K26").Value = 20000
$ws.Range("L26").Value = 26637.666
$ws.Range("M26").Value = -19720
$ws.Range("N26").Value = -27197.666
$ws.Range("H50").Value = 25689.428
$ws.Range("I50").Value = 20000
$ws.Range("J50").Value = 26637.666
$ws.Range("K50").Value = 20000
$ws.Range("L50").Value = 26637.666
$ws.Range("M50").Value = -19502
$ws.Range("N50").Value = -27633.666
$ws.Range("H58").Value = 40000
$ws.Range("J58").Value = 40000
$ws.Range("L58").Value = 40000
$ws.Range("N58").Value = -40554

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2551.2632
$ws.Range("I7").Value = 2608.6365
$ws.Range("J7").Value = 2472.375
$ws.Range("K7").Value = 2608.6365
$ws.Range("L7").Value = 2472.375
$ws.Range("M7").Value = -2496.6365
$ws.Range("N7").Value = -2696.375
$ws.Range("H40").Value = 2771.8572
$ws.Range("I40").Value = 2683.8333
$ws.Range("J40").Value = 3300
$ws.Range("K40").Value = 2683.8333
$ws.Range("L40").Value = 3300
$ws.Range("M40").Value = -2547.8333
$ws.Range("N40").Value = -3572
$ws.Range("H122").Value = 2392.3704
$ws.Range("I122").Value = 2209.0454
$ws.Range("J122").Value = 3199
$ws.Range("K122").Value = 6627.1362
$ws.Range("L122").Value = 9597
$ws.Range("M122").Value = -4177.1362
$ws.Range("N122").Value = -14497
$ws.Range("H126").Value = 2551.2632
$ws.Range("I126").Value = 2608.6365
$ws.Range("J126").Value = 2472.375
$ws.Range("K126").Value = 7825.9095
$ws.Range("L126").Value = 7417.125
$ws.Range("M126").Value = -5355.9095
$ws.Range("N126").Value = -12357.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2615.1428
$ws.Range("I96").Value = 2168.6667
$ws.Range("J96").Value = 2950
$ws.Range("K96").Value = 2168.6667
$ws.Range("L96").Value = 2950
$ws.Range("M96").Value = -795.6667000000002
$ws.Range("N96").Value = -5696
$ws.Range("H113").Value = 352.54285
$ws.Range("I113").Value = 236.61111
$ws.Range("K113").Value = 709.8333299999999
$ws.Range("M113").Value = 1460.16667
